$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 5942
$ws.Range("F6").Value = 2994
$ws.Range("F7").Value = 1282
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = 737
$ws.Range("F13").Value = 284
$ws.Range("F14").Value = 4378
$ws.Range("F15").Value = 4378
$ws.Range("F17").Value = 90
$ws.Range("F21").Value = 71
$ws.Range("F22").Value = 6729
$ws.Range("F23").Value = 234
$ws.Range("F24").Value = 106
$ws.Range("F25").Value = 464
$ws.Range("F26").Value = 1251
$ws.Range("F27").Value = 6257
$ws.Range("F30").Value = 1866
$ws.Range("F31").Value = 6003
$ws.Range("F36").Value = 426
$ws.Range("F37").Value = 4439
$ws.Range("F39").Value = 189
$ws.Range("F41").Value = 18
$ws.Range("F42").Value = 4
$ws.Range("F43").Value = 2412
$ws.Range("F44").Value = 25
$ws.Range("F45").Value = 46
$ws.Range("F48").Value = 340
$ws.Range("F49").Value = 2065
$ws.Range("F50").Value = 6

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 31
$ws.Range("F9").Value = 4

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 5942
$ws.Range("F5").Value = 2994
$ws.Range("F6").Value = 1282
$ws.Range("F11").Value = 31
$ws.Range("F12").Value = 284
$ws.Range("F13").Value = 4378
$ws.Range("F14").Value = 4378
$ws.Range("F16").Value = 90
$ws.Range("F20").Value = 71
$ws.Range("F21").Value = 6729
$ws.Range("F22").Value = 234
$ws.Range("F23").Value = 106
$ws.Range("F24").Value = 464
$ws.Range("F25").Value = 1251
$ws.Range("F27").Value = 6257
$ws.Range("F29").Value = 1866
$ws.Range("F31").Value = 6003
$ws.Range("F33").Value = 4
$ws.Range("F37").Value = 426
$ws.Range("F38").Value = 4440
$ws.Range("F40").Value = 189
$ws.Range("F42").Value = 18
$ws.Range("F44").Value = 2412
$ws.Range("F45").Value = 25
$ws.Range("F46").Value = 46
$ws.Range("F49").Value = 340
